$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '62.800.88'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +1.31%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.468.66'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +1.19%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '414.02'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +1.16%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '130.24'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +1.23%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.623'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -1.46%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.725'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -1.76%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.146'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +3.45%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '42.59'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -0.45%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '9.69'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  +5.96%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.141'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -0.12%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '20.48'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -4.41%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '3.449.20'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -0.35%  '
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +1.87%  '
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.93%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '62.663.64'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +1.16%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '466.40'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  +2.92%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '90.78'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -1.27%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.27'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +1.84%  '
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +2.38%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '10.57'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +19.96%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '3.30'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  +2.14%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '33.40'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +1.38%  '
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +0.35%  '
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -1.00%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '11.99'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  -0.09%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '2.65'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -2.55%  '
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  -2.44%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '40.81'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -5.33%  '
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +0.12%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '58.57'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  +7.60%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.0489'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  -2.67%  '
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -1.11%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '146.13'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  +2.76%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.08'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +4.12%  '
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +1.98%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.40'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  +12.18%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.0₃0556'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +31.30%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '16.34'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -2.06%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '22.18'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -0.48%  '
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +1.30%  '
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000219'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -0.99%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.009.92'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +1.19%  '
$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.09'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  +5.51%  '
$ws.Cells.Item(39, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.998'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  +0.03%  '
$ws.Cells.Item(41, 2).Value = 'TheGraph'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.321'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -0.71%  '
$ws.Cells.Item(42, 2).Value = 'Stellar'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.134'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -0.84%  '
$ws.Cells.Item(43, 2).Value = 'WEMIXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.69'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  +7.34%  '
